# The document's "first page" footer, "default" footer and "first page"
# header each contain a single inline picture (the Pearson / BTEC logos).
# The edit simply renames those pictures:
#   footer (first)   : image2.png -> image1.png
#   footer (default) : image2.png -> image1.png
#   header (first)   : image1.jpg -> image2.jpg
#
# InlineShape has no writable "Name" property (same as real Word), so each
# picture is temporarily converted to a floating Shape - which does expose
# .Name - renamed, then converted back to an InlineShape so the drawing
# stays wp:inline (matching the original layout). AlternativeText (the
# "descr" attribute) is preserved explicitly because ConvertToShape does
# not carry it over on its own.

function Rename-InlinePicture($range, $newName) {
    $shape = $range.InlineShapes.Item(1)
    $savedAltText = $shape.AlternativeText

    $floating = $shape.ConvertToShape()
    $floating.Name = $newName
    $floating.AlternativeText = $savedAltText

    [void]$floating.ConvertToInlineShape()
}

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
Rename-InlinePicture $section.Footers.Item(2).Range "image1.png"
Rename-InlinePicture $section.Footers.Item(1).Range "image1.png"
Rename-InlinePicture $section.Headers.Item(2).Range "image2.jpg"
